# Helper: classic VBA-style RGB() packer (r + g*256 + b*65536)
function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 5's table: switch from the custom "Table_0" style to the built-in
#    table style {D4B49255-DE0D-4C5C-AFEE-07B252013A08}.
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{D4B49255-DE0D-4C5C-AFEE-07B252013A08}")
    }
}

# ---------------------------------------------------------------------------
# 2) Re-colour the deck's theme (slide master's theme, theme1.xml) from the
#    "Integral / Red Violet" palette to the default "Office" palette - i.e.
#    switch the presentation Design to the stock Office Theme colours.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$tcs = $slide.ThemeColorScheme

$tcs.Item(1).RGB  = RGB 0x00 0x00 0x00   # Dark 1
$tcs.Item(2).RGB  = RGB 0xFF 0xFF 0xFF   # Light 1
$tcs.Item(3).RGB  = RGB 0x44 0x54 0x6A   # Dark 2
$tcs.Item(4).RGB  = RGB 0xE7 0xE6 0xE6   # Light 2
$tcs.Item(5).RGB  = RGB 0x5B 0x9B 0xD5   # Accent 1
$tcs.Item(6).RGB  = RGB 0xED 0x7D 0x31   # Accent 2
$tcs.Item(7).RGB  = RGB 0xA5 0xA5 0xA5   # Accent 3
$tcs.Item(8).RGB  = RGB 0xFF 0xC0 0x00   # Accent 4
$tcs.Item(9).RGB  = RGB 0x44 0x72 0xC4   # Accent 5
$tcs.Item(10).RGB = RGB 0x70 0xAD 0x47   # Accent 6
$tcs.Item(11).RGB = RGB 0x05 0x63 0xC1   # Hyperlink
$tcs.Item(12).RGB = RGB 0x95 0x4F 0x72   # Followed hyperlink
